# C5-PowerPoint.pptx edit
# 1) Re-style the table on slide 6 with the new table style GUID.
# 2) Swap the theme color palettes: the deck's (single, shared) theme
#    currently carries the "Integral" palette; it becomes the "Office
#    Theme" palette (this is the palette the slide master/visible slides
#    actually render with).

function RGBVal([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

function Set-ThemeColors($theme, [string[]]$hexColors) {
    $tcs = $theme.ThemeColorScheme
    for ($i = 0; $i -lt $hexColors.Count; $i++) {
        $tcs.Item($i + 1).RGB = RGBVal $hexColors[$i]
    }
}

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 ---------------------------------------
$slide6 = $p.Slides.Item(6)
$tableShape = $slide6.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{6FA941D9-5E25-4C36-B342-BF090F36C9CD}")

# --- 2. Swap the theme color palette -----------------------------------
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    "000000", "FFFFFF", "44546A", "E7E6E6",
    "5B9BD5", "ED7D31", "A5A5A5", "FFC000", "4472C4", "70AD47",
    "0563C1", "954F72"
)

# The presentation's theme (used by the slide master, notes master and
# handout master alike) currently carries the Integral palette -> Office
# Theme palette.
Set-ThemeColors $p.SlideMaster.Theme $officeColors
